{"js": "// Replace each unique old text with its corresponding new text.\n// Every \"old\" string below occurs exactly once in the document body,\n// so a simple search + insertText(\"Replace\") round-trip is safe and\n// preserves the existing run formatting (fonts, size, etc.).\nconst replacements = [\n  [\"2025-01-19 Sunday\", \"2025-01-20 Monday\"],\n  [\"183\u00f73=61, 0\", \"246\u00f79=27, 3\"],\n  [\"311\u00f77=44, 3\", \"264\u00f77=37, 5\"],\n  [\"595\u00f77=85, 0\", \"417\u00f74=104, 1\"],\n  [\"446\u00f75=89, 1\", \"735\u00f78=91, 7\"],\n  [\"793\u00f76=132, 1\", \"501\u00f77=71, 4\"],\n  [\"592\u00f75=118, 2\", \"793\u00f75=158, 3\"],\n  [\"315\u00f74=78, 3\", \"813\u00f74=203, 1\"],\n  [\"582\u00f77=83, 1\", \"892\u00f77=127, 3\"],\n  [\"849\u00f72=424, 1\", \"303\u00f74=75, 3\"],\n  [\"880\u00f72=440, 0\", \"377\u00f78=47, 1\"],\n  [\"703\u00f75=140, 3\", \"494\u00f78=61, 6\"],\n  [\"753\u00f78=94, 1\", \"526\u00f77=75, 1\"],\n  [\"115\u00f74=28, 3\", \"673\u00f72=336, 1\"],\n  [\"614\u00f73=204, 2\", \"138\u00f74=34, 2\"],\n  [\"821\u00f77=117, 2\", \"554\u00f72=277, 0\"],\n  [\"865\u00f76=144, 1\", \"807\u00f73=269, 0\"],\n  [\"613\u00f76=102, 1\", \"797\u00f77=113, 6\"],\n  [\"254\u00f73=84, 2\", \"166\u00f75=33, 1\"],\n  [\"890\u00f77=127, 1\", \"857\u00f77=122, 3\"],\n  [\"732\u00f73=244, 0\", \"379\u00f75=75, 4\"],\n  [\"731\u00f77=104, 3\", \"333\u00f75=66, 3\"],\n  [\"515\u00f72=257, 1\", \"721\u00f74=180, 1\"],\n  [\"606\u00f75=121, 1\", \"370\u00f73=123, 1\"],\n  [\"714\u00f78=89, 2\", \"753\u00f77=107, 4\"],\n  [\"463\u00f77=66, 1\", \"108\u00f72=54, 0\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for: \" + oldText + \" (found \" + results.items.length + \")\");\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each pair is a unique old string in the document paired with its\n# replacement text, taken from the commit diff.\n$replacements = @(\n    @('2025-01-19 Sunday', '2025-01-20 Monday'),\n    @('183\u00f73=61, 0', '246\u00f79=27, 3'),\n    @('311\u00f77=44, 3', '264\u00f77=37, 5'),\n    @('595\u00f77=85, 0', '417\u00f74=104, 1'),\n    @('446\u00f75=89, 1', '735\u00f78=91, 7'),\n    @('793\u00f76=132, 1', '501\u00f77=71, 4'),\n    @('592\u00f75=118, 2', '793\u00f75=158, 3'),\n    @('315\u00f74=78, 3', '813\u00f74=203, 1'),\n    @('582\u00f77=83, 1', '892\u00f77=127, 3'),\n    @('849\u00f72=424, 1', '303\u00f74=75, 3'),\n    @('880\u00f72=440, 0', '377\u00f78=47, 1'),\n    @('703\u00f75=140, 3', '494\u00f78=61, 6'),\n    @('753\u00f78=94, 1', '526\u00f77=75, 1'),\n    @('115\u00f74=28, 3', '673\u00f72=336, 1'),\n    @('614\u00f73=204, 2', '138\u00f74=34, 2'),\n    @('821\u00f77=117, 2', '554\u00f72=277, 0'),\n    @('865\u00f76=144, 1', '807\u00f73=269, 0'),\n    @('613\u00f76=102, 1', '797\u00f77=113, 6'),\n    @('254\u00f73=84, 2', '166\u00f75=33, 1'),\n    @('890\u00f77=127, 1', '857\u00f77=122, 3'),\n    @('732\u00f73=244, 0', '379\u00f75=75, 4'),\n    @('731\u00f77=104, 3', '333\u00f75=66, 3'),\n    @('515\u00f72=257, 1', '721\u00f74=180, 1'),\n    @('606\u00f75=121, 1', '370\u00f73=123, 1'),\n    @('714\u00f78=89, 2', '753\u00f77=107, 4'),\n    @('463\u00f77=66, 1', '108\u00f72=54, 0')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $ok = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $oldText\"\n    }\n}\n\n"}
